# Dataframe ST.xlsx update:
#  - Sheet3 lookup table (B20:B36) gets refreshed figures for the new day.
#  - Sheet1 VLOOKUP columns (CB/CC) recalc automatically from that change.
#  - A new snapshot column ("14-nov") is appended after CN, freezing the
#    freshly recalculated values (mirrors how CN itself froze "13-nov").

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet3 = $wb.Worksheets.Item("Sheet3")

# --- 1. Update the raw lookup values on Sheet3 (rows 20-36, column B) ---
$sheet3Updates = @{
    20 = 12.227808641918193
    23 = 7.8693392877096953
    24 = 6.9667987848559561
    26 = 6.8454685958566852
    27 = 4.9110207616356449
    28 = 6.7020066671549401
    29 = 6.2287312161842721
    31 = 1.08051533333328
    32 = 16.390779966015895
    33 = 12.486550123753087
    34 = 7.1286032676410462
    35 = 7.0605404447186855
    36 = 44.650810748545069
}

foreach ($row in $sheet3Updates.Keys) {
    $sheet3.Cells.Item($row, 2).Value = $sheet3Updates[$row]
}

# Force a full recalculation so Sheet1's VLOOKUP formulas (CB/CC) pick up
# the new Sheet3 figures before we snapshot them into the new column.
$excel.CalculateFullRebuild()

# --- 2. Append the new "14-nov" snapshot column (CO) on Sheet1 ---
$sheet1.Cells.Item(1, 93).Value = "14-nov"

for ($row = 2; $row -le 18; $row++) {
    $recalced = $sheet1.Cells.Item($row, 80).Value()
    $target = $sheet1.Cells.Item($row, 93)
    $target.Value = $recalced
    $target.NumberFormat = $sheet1.Cells.Item($row, 92).NumberFormat
}

# --- 3. Cosmetic: leave the selection where the editor last clicked ---
$sheet1.Range("CJ26").Select()
